# Auto-generated edit script: appends newly logged sensor rows to
# the ALERTS, PIR, Humidity, and Temperature sheets.
$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: rows 7-9 ---
$ws = $wb.Worksheets.Item("ALERTS")
$newRows = @(
    @("2026-02-06", "09:51:35", "09:00", "Bathroom", "MINIMAL", "MINIMAL ALERT: Bathroom occupied, no motion > 20s."),
    @("2026-02-06", "09:51:55", "09:00", "Bathroom", "MODERATE", "MODERATE ALERT: Bathroom occupied, no motion > 40s."),
    @("2026-02-06", "09:52:15", "09:00", "Bathroom", "CRITICAL", "CRITICAL ALERT: Bathroom occupied, no motion > 60s.")
)
$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# --- PIR sheet: rows 138-150 ---
$ws = $wb.Worksheets.Item("PIR")
$newRows = @(
    @("2026-02-06", "09:51:33", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:34", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:38", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:43", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:48", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:53", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:51:58", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:03", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:08", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:13", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:18", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:23", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:52:28", "09:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 138
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# --- Humidity sheet: rows 69-79 ---
$ws = $wb.Worksheets.Item("Humidity")
$newRows = @(
    @("2026-02-06", "09:51:33", "09:00", "Bathroom", "75.3%", "Active"),
    @("2026-02-06", "09:51:36", "09:00", "Bathroom", "74.1%", "Active"),
    @("2026-02-06", "09:51:41", "09:00", "Bathroom", "73.3%", "Active"),
    @("2026-02-06", "09:51:46", "09:00", "Bathroom", "72.9%", "Active"),
    @("2026-02-06", "09:51:56", "09:00", "Bathroom", "70.6%", "Active"),
    @("2026-02-06", "09:52:01", "09:00", "Bathroom", "71.8%", "Active"),
    @("2026-02-06", "09:52:06", "09:00", "Bathroom", "71.8%", "Active"),
    @("2026-02-06", "09:52:16", "09:00", "Bathroom", "70.6%", "Active"),
    @("2026-02-06", "09:52:21", "09:00", "Bathroom", "71.6%", "Active"),
    @("2026-02-06", "09:52:26", "09:00", "Bathroom", "70.7%", "Active"),
    @("2026-02-06", "09:52:31", "09:00", "Bathroom", "71.5%", "Active")
)
$startRow = 69
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# --- Temperature sheet: rows 69-79 ---
$ws = $wb.Worksheets.Item("Temperature")
$newRows = @(
    @("2026-02-06", "09:51:33", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:51:36", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:51:41", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:51:47", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:51:56", "09:00", "Bathroom", "27.6C", "Active"),
    @("2026-02-06", "09:52:01", "09:00", "Bathroom", "27.6C", "Active"),
    @("2026-02-06", "09:52:06", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:52:16", "09:00", "Bathroom", "27.6C", "Active"),
    @("2026-02-06", "09:52:22", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:52:27", "09:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "09:52:32", "09:00", "Bathroom", "27.7C", "Active")
)
$startRow = 69
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

Write-Host "Appended new sensor log rows to ALERTS, PIR, Humidity, and Temperature sheets."
